# Updated cryptos list on Sat Jul 20 02:41:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.597.44"
$ws.Range("E2").Value = "  +4.42%  "
$ws.Range("D3").Value = "3.493.45"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'588.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.24%  "
$ws.Range("D6").Value = "'171.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.75%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.491.92"
$ws.Range("E8").Value = "  +2.31%  "
$ws.Range("D9").Value = "'0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  +4.92%  "
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").Value = "4.093.37"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "'28.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("D16").Value = "66.616.51"
$ws.Range("E16").Value = "  +4.30%  "
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "3.486.87"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").Value = "'6.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.39%  "
$ws.Range("D20").Value = "'13.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").Value = "'387.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("D22").Value = "'7.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").Value = "'73.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("E26").Value = "  +5.10%  "
$ws.Range("D27").Value = "'10.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.41%  "
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'6.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.26%  "
$ws.Range("D31").Value = "'1.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.16%  "
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("E34").Value = "  +5.60%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'1.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.66%  "
$ws.Range("D37").Value = "'163.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("E38").Value = "  +5.48%  "
$ws.Range("E39").Value = "  +5.11%  "
$ws.Range("D40").Value = "'0.0745"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +5.16%  "
$ws.Range("D42").Value = "'26.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "'6.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Value = "2.805.25"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "'26.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.92%  "
$ws.Range("D46").Value = "'42.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "'356.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.17%  "
$ws.Range("E50").Value = "  +3.92%  "
$ws.Range("D51").Value = "'33.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.45%  "

# Rows 48 and 49 swap places (VeChain <-> dogwifhat) with updated price/volume
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.28%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0309"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.45%  "

